$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "아미야는 지마 일행 중 유일하게 후방 지원 부서에서 일하던 나탈리야의 전투원 전환 신청을 받게 된다.`n"
$ws.Range("D3").Value = "이에 따라 아미야는 박사에게 이 신청자를 직접 만나봐 달라 부탁하게 되는데……`n"
